$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 229 - this pushes the existing rows 229:263
# down to 230:264 (and extends the sheet dimension to A1:R264).
$ws.Rows.Item(229).Insert()

# Populate the newly inserted row 229 with the new weekly record
# (same data as the row that is now 231, i.e. the previous row 229's
# neighbour, but with an updated date).
$ws.Cells.Item(229, 1).Value = 9
$ws.Cells.Item(229, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(229, 3).Value = "Metropolitana"
$ws.Cells.Item(229, 4).Value = 44474
$ws.Cells.Item(229, 5).Value = 13
$ws.Cells.Item(229, 6).Value = 100112039
$ws.Cells.Item(229, 7).Value = "Ciboulette"
$ws.Cells.Item(229, 8).Value = "Sin especificar"
$ws.Cells.Item(229, 9).Value = "Primera"
$ws.Cells.Item(229, 10).Value = 250
$ws.Cells.Item(229, 11).Value = 1000
$ws.Cells.Item(229, 12).Value = 1200
$ws.Cells.Item(229, 13).Value = 1100
$ws.Cells.Item(229, 14).Value = "$/docena de atados"
$ws.Cells.Item(229, 15).Value = "Región Metropolitana"
$ws.Cells.Item(229, 16).Value = 367
$ws.Cells.Item(229, 17).Value = 3
$ws.Cells.Item(229, 18).Value = "Hortaliza"

# Match the date-style used by the rest of column D.
$ws.Cells.Item(229, 4).NumberFormat = $ws.Cells.Item(230, 4).NumberFormat
